# Auto-generated edit script: update crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.326.03"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.570.30"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.61"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.77"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.559.23"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +8.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.82"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000307"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.56"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.139.84"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.48"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.322.03"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.560.43"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "544.94"
$ws.Range("E21").Value = "  +11.38%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.03"
$ws.Range("E23").Value = "  -7.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.68"
$ws.Range("E24").Value = "  +8.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.94"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.95"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +4.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.00"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.26"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.55"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.24"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "550.98"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.21"
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.416"
$ws.Range("E37").Value = "  +5.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.44"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0770"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.372.01"
$ws.Range("E42").Value = "  +3.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.08"
$ws.Range("E44").Value = "  -6.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.58"
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0447"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.21"
$ws.Range("E48").Value = "  -4.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.137"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.44"
$ws.Range("E51").Value = "  +18.84%  "
